# Backup QR Scanner data - 2025-12-16T08:52:45.572Z - Cache Bust: 1765875165572
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Session" to "Neurology"
$ws.Name = "Neurology"

# Append the new scanner log rows (37-40)
$rows = @(
    @("190333", "Neurology", "16/12/2025", "10:13:46", "Manual", "emp17.farah.a.youssef@gmail.com"),
    @("191007", "Neurology", "16/12/2025", "10:16:24", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("202051", "Neurology", "16/12/2025", "10:48:43", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("210728", "Neurology", "16/12/2025", "10:52:36", "Scan",   "emp17.farah.a.youssef@gmail.com")
)

$styleRef = $ws.Cells.Item(2, 1)

$startRow = 37
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $value = $data[$c - 1]
        # Force numeric-looking values (e.g. student IDs) to stay text, as in
        # the rest of the sheet, by using a leading apostrophe, then reset the
        # cell's style to match the existing data rows (avoids picking up the
        # quote-prefix style flag).
        $cell.Value = "'" + $value
        $cell.Style = $styleRef.Style
    }
}
